# Updated cryptos list on Wed Dec 20 13:32:54 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) hold plain text (not numbers), so a
# leading apostrophe forces text interpretation for number-like values; the
# cell style is then reset to Normal so no stray "quote prefix" style index
# is introduced (the source cells carry no explicit style).

# Row 2
$ws.Range('D2').Value = "'43.389.52"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.99%  "
$ws.Range('E2').Style = 'Normal'
# Row 3
$ws.Range('D3').Value = "'2.233.86"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.08%  "
$ws.Range('E3').Style = 'Normal'
# Row 4
$ws.Range('E4').Value = "'  -0.17%  "
$ws.Range('E4').Style = 'Normal'
# Row 5
$ws.Range('D5').Value = "'258.50"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +2.46%  "
$ws.Range('E5').Style = 'Normal'
# Row 6
$ws.Range('D6').Value = "'0.625"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +1.49%  "
$ws.Range('E6').Style = 'Normal'
# Row 7
$ws.Range('D7').Value = "'78.47"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +4.58%  "
$ws.Range('E7').Style = 'Normal'
# Row 8
$ws.Range('E8').Value = "'  -0.11%  "
$ws.Range('E8').Style = 'Normal'
# Row 9
$ws.Range('D9').Value = "'0.601"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +1.00%  "
$ws.Range('E9').Style = 'Normal'
# Row 10
$ws.Range('D10').Value = "'43.16"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +4.64%  "
$ws.Range('E10').Style = 'Normal'
# Row 11
$ws.Range('D11').Value = "'0.0925"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.22%  "
$ws.Range('E11').Style = 'Normal'
# Row 12
$ws.Range('D12').Value = "'7.11"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +3.30%  "
$ws.Range('E12').Style = 'Normal'
# Row 13
$ws.Range('E13').Value = "'  +1.41%  "
$ws.Range('E13').Style = 'Normal'
# Row 14
$ws.Range('D14').Value = "'2.564.71"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.14%  "
$ws.Range('E14').Style = 'Normal'
# Row 15
$ws.Range('D15').Value = "'14.67"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +1.11%  "
$ws.Range('E15').Style = 'Normal'
# Row 16
$ws.Range('D16').Value = "'2.237.77"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.09%  "
$ws.Range('E16').Style = 'Normal'
# Row 17
$ws.Range('D17').Value = "'0.797"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +1.02%  "
$ws.Range('E17').Style = 'Normal'
# Row 18
$ws.Range('D18').Value = "'43.278.82"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.91%  "
$ws.Range('E18').Style = 'Normal'
# Row 19
$ws.Range('D19').Value = "'0.0000105"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.94%  "
$ws.Range('E19').Style = 'Normal'
# Row 20
$ws.Range('D20').Value = "'71.38"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.23%  "
$ws.Range('E20').Style = 'Normal'
# Row 21
$ws.Range('D21').Value = "'6.04"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +1.46%  "
$ws.Range('E21').Style = 'Normal'
# Row 22
$ws.Range('E22').Value = "'  +5.65%  "
$ws.Range('E22').Style = 'Normal'
# Row 23
$ws.Range('D23').Value = "'231.81"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +0.87%  "
$ws.Range('E23').Style = 'Normal'
# Row 24
$ws.Range('D24').Value = "'9.35"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +0.03%  "
$ws.Range('E24').Style = 'Normal'
# Row 25
$ws.Range('E25').Value = "'  -0.12%  "
$ws.Range('E25').Style = 'Normal'
# Row 26
$ws.Range('D26').Value = "'42.19"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +8.11%  "
$ws.Range('E26').Style = 'Normal'
# Row 27
$ws.Range('D27').Value = "'10.85"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +1.34%  "
$ws.Range('E27').Style = 'Normal'
# Row 28
$ws.Range('E28').Value = "'  -2.20%  "
$ws.Range('E28').Style = 'Normal'
# Row 29
$ws.Range('D29').Value = "'2.22"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.08%  "
$ws.Range('E29').Style = 'Normal'
# Row 30
$ws.Range('E30').Value = "'  -1.34%  "
$ws.Range('E30').Style = 'Normal'
# Row 31
$ws.Range('D31').Value = "'173.53"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.93%  "
$ws.Range('E31').Style = 'Normal'
# Row 32
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = "'0.0871"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +9.66%  "
$ws.Range('E32').Style = 'Normal'
# Row 33
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = "'20.23"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +0.28%  "
$ws.Range('E33').Style = 'Normal'
# Row 34
$ws.Range('E34').Value = "'  +0.51%  "
$ws.Range('E34').Style = 'Normal'
# Row 36
$ws.Range('D36').Value = "'0.0370"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +13.36%  "
$ws.Range('E36').Style = 'Normal'
# Row 37
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = "'4.46"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +0.12%  "
$ws.Range('E37').Style = 'Normal'
# Row 38
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = "'0.108"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -5.12%  "
$ws.Range('E38').Style = 'Normal'
# Row 39
$ws.Range('D39').Value = "'13.33"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +7.84%  "
$ws.Range('E39').Style = 'Normal'
# Row 40
$ws.Range('D40').Value = "'2.93"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +20.45%  "
$ws.Range('E40').Style = 'Normal'
# Row 41
$ws.Range('E41').Value = "'  +1.92%  "
$ws.Range('E41').Style = 'Normal'
# Row 42
$ws.Range('E42').Value = "'  +0.04%  "
$ws.Range('E42').Style = 'Normal'
# Row 43
$ws.Range('D43').Value = "'61.78"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +3.58%  "
$ws.Range('E43').Style = 'Normal'
# Row 44
$ws.Range('D44').Value = "'5.35"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -0.09%  "
$ws.Range('E44').Style = 'Normal'
# Row 45
$ws.Range('D45').Value = "'103.88"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +0.73%  "
$ws.Range('E45').Style = 'Normal'
# Row 46
$ws.Range('D46').Value = "'8.63"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.32%  "
$ws.Range('E46').Style = 'Normal'
# Row 47
$ws.Range('D47').Value = "'0.474"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -2.79%  "
$ws.Range('E47').Style = 'Normal'
# Row 48
$ws.Range('D48').Value = "'0.0982"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -0.29%  "
$ws.Range('E48').Style = 'Normal'
# Row 49
$ws.Range('E49').Value = "'  +0.64%  "
$ws.Range('E49').Style = 'Normal'
# Row 50
$ws.Range('E50').Value = "'  +1.02%  "
$ws.Range('E50').Style = 'Normal'
# Row 51
$ws.Range('E51').Value = "'  +22.90%  "
$ws.Range('E51').Style = 'Normal'
